# Refresh cryptocurrency market data (prices / 1h volume %) pulled from coinranking.com
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.737.77"
$ws.Range("E2").Value = "  +3.72%  "
$ws.Range("D3").Value = "3.007.93"
$ws.Range("E3").Value = "  +2.42%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.69"
$ws.Range("E5").Value = "  +2.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.39"
$ws.Range("E6").Value = "  +7.30%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.520"
$ws.Range("E8").Value = "  +1.45%  "
$ws.Range("D9").Value = "2.995.91"
$ws.Range("E9").Value = "  +2.27%  "
$ws.Range("E10").Value = "  +6.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.35"
$ws.Range("E11").Value = "  +12.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000231"
$ws.Range("E13").Value = "  +4.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.93"
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("D16").Value = "3.505.09"
$ws.Range("E16").Value = "  +2.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.21"
$ws.Range("E17").Value = "  +4.59%  "
$ws.Range("D18").Value = "3.009.32"
$ws.Range("E18").Value = "  +2.61%  "
$ws.Range("D19").Value = "59.716.60"
$ws.Range("E19").Value = "  +3.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "436.97"
$ws.Range("E20").Value = "  +4.70%  "
$ws.Range("E21").Value = "  +2.93%  "
$ws.Range("E22").Value = "  +4.27%  "
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.36"
$ws.Range("E24").Value = "  +2.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.57"
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  +10.96%  "
$ws.Range("E29").Value = "  +3.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.81"
$ws.Range("E30").Value = "  +5.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.34"
$ws.Range("E31").Value = "  +5.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.95"
$ws.Range("E32").Value = "  +2.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.105"
$ws.Range("E33").Value = "  +8.80%  "
$ws.Range("E34").Value = "  +14.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  +6.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.92"
$ws.Range("E36").Value = "  +4.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.10"
$ws.Range("E37").Value = "  +1.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.19"
$ws.Range("E38").Value = "  +2.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.57"
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("E40").Value = "  +10.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "402.47"
$ws.Range("E41").Value = "  +6.90%  "
$ws.Range("E42").Value = "  +2.22%  "
$ws.Range("D43").Value = "2.764.13"
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("E45").Value = "  +6.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "123.50"
$ws.Range("E47").Value = "  +0.88%  "
$ws.Range("E48").Value = "  +3.93%  "
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.08"
$ws.Range("E49").Value = "  +20.31%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.110"
$ws.Range("E50").Value = "  +1.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.62"
$ws.Range("E51").Value = "  +2.51%  "
